$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I mirrors the width/formatting of column E
$ws.Columns.Item(9).ColumnWidth = 10.43

# I5: a date value (2024-01-15), formatted as a date (reuses the existing date style)
$ws.Range("I5").Value = 45306
$ws.Range("I5").NumberFormat = "m/d/yy"

# I6: formula referencing I5, building a URL string
$ws.Range("I6").Formula = '="http://124.55.158.229:4752/get_main_data/" & TEXT(I5, "yyyy-mm-dd") & "/KRW-USD"'

# Update the active selection to I9, matching the diff
$ws.Range("I9").Select()
